$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.252.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "'3.369.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'590.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.82%  "
$ws.Range("D6").Value = "'188.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").Value = "'0.603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.43%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.186"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'47.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").Value = "'655.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.31%  "
$ws.Range("D14").Value = "'3.906.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "'8.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "'67.307.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "'18.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'3.369.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "'11.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").Value = "'0.911"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").Value = "'18.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "'101.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'4.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").Value = "'2.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").Value = "'9.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").Value = "'32.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.52%  "
$ws.Range("D29").Value = "'8.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").Value = "'6.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.33%  "
$ws.Range("D31").Value = "'617.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.19%  "
$ws.Range("D32").Value = "'3.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'3.892.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "'55.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "'2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.85%  "
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").Value = "'33.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "'0.0₃0713"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").Value = "'0.348"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("D44").Value = "'3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("D46").Value = "'0.130"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").Value = "'2.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -16.39%  "
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.26%  "
$ws.Range("D51").Value = "'130.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.88%  "
